# "Added two new Mac-Addresses"
# Appends two new data rows (31 and 32) to the reg_center_user_h test-data
# sheet, following the exact same pattern as the existing rows: regcntr_id
# 10001, the next sequential usr_id values (110030 / 110031), lang_code
# "eng", is_active TRUE, cr_by "superadmin", and cr_dtimes/eff_dtimes both
# "now()".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 - new Mac-Address record
$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 110030
$ws.Cells.Item(31, 3).Value = "eng"
$ws.Cells.Item(31, 4).Value = $true
$ws.Cells.Item(31, 5).Value = "superadmin"
$ws.Cells.Item(31, 6).Value = "now()"
$ws.Cells.Item(31, 7).Value = "now()"

# Row 32 - new Mac-Address record
$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 110031
$ws.Cells.Item(32, 3).Value = "eng"
$ws.Cells.Item(32, 4).Value = $true
$ws.Cells.Item(32, 5).Value = "superadmin"
$ws.Cells.Item(32, 6).Value = "now()"
$ws.Cells.Item(32, 7).Value = "now()"

# Mirror the view state left after the edit: the window had scrolled down
# so row 19 is the top visible row, with E28 as the active selected cell.
try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Window scroll position isn't always settable in every host; ignore.
}
$ws.Range("E28").Select()
